$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.008.21'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '2.643.59'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("E5").Value = '  +1.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = '2.656.00'
$ws.Range("E9").Value = '  -0.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.33'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.105'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.336'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.08%  '
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").Value = '3.102.71'
$ws.Range("E14").Value = '  -0.19%  '
$ws.Range("D15").Value = '59.632.20'
$ws.Range("E15").Value = '  +1.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.95'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.85%  '
$ws.Range("E17").Value = '  -0.59%  '
$ws.Range("D18").Value = '2.640.29'
$ws.Range("E18").Value = '  -1.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '350.80'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.63%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.48'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.63%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.34'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.19'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.58%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.416'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.164'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.54%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.995'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.05%  '
$ws.Range("D28").Value = '0.0₃0808'
$ws.Range("E28").Value = '  -2.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.09'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.32'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.11%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.58'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.970'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.04'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.16'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.57%  '
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.852'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.50%  '
$ws.Range("B39").Value = 'OKB'
$ws.Range("C39").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.62'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("E40").Value = '  +0.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.65'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '281.53'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0987'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '19.82'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.43%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.603'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.36%  '
$ws.Range("D47").Value = '2.090.25'
$ws.Range("E47").Value = '  +5.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0528'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.67%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0231'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '10.30'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.70'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.16%  '
